# Apply commit "#5: fund, bonds, otherbonds, antique done" to the workbook.
# Sheet7 = 基金受益憑證 (fund), Sheet8 = 具有相當價值之財產 (otherbonds), Sheet9 = 保險.
# All three sheets gain the full set of "normal" metadata columns
# (property_category/category/date/legislator_name/legislator_id/source_file/index)
# that the other sheets (土地/建物/汽車/存款/股票/債券) already have, and sheet7
# additionally gains explicit currency/total columns (shifting the old F/G data
# right by one) plus a numeric face_value placeholder column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet7: 基金受益憑證 (fund)
# ---------------------------------------------------------------------------
$ws7 = $wb.Worksheets.Item(7)

# Fix header row (B1:G1 previously held stray data values instead of labels)
# and extend it with the new trailing metadata headers (H1:O1).
$headers7 = @("name","owner","dealer","quantity","face_value","currency","total","property_category","category","date","legislator_name","legislator_id","source_file","index")
for ($i = 0; $i -lt $headers7.Length; $i++) {
    $col = 2 + $i
    $cell = $ws7.Cells.Item(1, $col)
    if ($col -gt 7) {
        $ws7.Range("B1").Copy()
        $cell.PasteSpecial(-4122)
    }
    $cell.Value = $headers7[$i]
}

# Data rows 2-8: shift currency (old F) -> G, total (old G) -> H, set F to 0,
# then append the metadata columns I-O.
$rows7 = @(117,118,119,120,121,122,123)
for ($r = 0; $r -lt $rows7.Length; $r++) {
    $row = 2 + $r
    $idx = $rows7[$r]

    $currency = $ws7.Cells.Item($row, 6).Value()
    $total = $ws7.Cells.Item($row, 7).Value()

    $ws7.Cells.Item($row, 6).Value = 0
    $ws7.Cells.Item($row, 7).Value = $currency
    $ws7.Cells.Item($row, 8).Value = $total

    for ($c = 9; $c -le 15; $c++) {
        $ws7.Range("G2").Copy()
        $ws7.Cells.Item($row, $c).PasteSpecial(-4122)
    }

    $ws7.Cells.Item($row, 9).Value = "fund"
    $ws7.Cells.Item($row, 10).Value = "normal"
    $ws7.Cells.Item($row, 11).Value = "2012-05-01"
    $ws7.Cells.Item($row, 12).Value = "黃偉哲"
    $ws7.Cells.Item($row, 13).Value = 1367
    $ws7.Cells.Item($row, 14).Value = "tmp62651"
    $ws7.Cells.Item($row, 15).Value = $idx
}

# ---------------------------------------------------------------------------
# Sheet8: 具有相當價值之財產 (otherbonds)
# ---------------------------------------------------------------------------
$ws8 = $wb.Worksheets.Item(8)

$headers8 = @("name","quantity","owner","total","property_category","category","date","legislator_name","legislator_id","source_file","index")
for ($i = 0; $i -lt $headers8.Length; $i++) {
    $col = 2 + $i
    $cell = $ws8.Cells.Item(1, $col)
    if ($col -gt 5) {
        $ws8.Range("B1").Copy()
        $cell.PasteSpecial(-4122)
    }
    $cell.Value = $headers8[$i]
}

for ($c = 6; $c -le 12; $c++) {
    $ws8.Range("B2").Copy()
    $ws8.Cells.Item(2, $c).PasteSpecial(-4122)
}
$ws8.Cells.Item(2, 6).Value = "otherbonds"
$ws8.Cells.Item(2, 7).Value = "normal"
$ws8.Cells.Item(2, 8).Value = "2012-05-01"
$ws8.Cells.Item(2, 9).Value = "黃偉哲"
$ws8.Cells.Item(2, 10).Value = 1367
$ws8.Cells.Item(2, 11).Value = "tmp62651"
$ws8.Cells.Item(2, 12).Value = 131

# ---------------------------------------------------------------------------
# Sheet9: 保險 - no structural change, values are unchanged (the shared
# string table simply grew), so nothing further to do here.
# ---------------------------------------------------------------------------
